$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "Report" to "Sheet1"
$ws.Name = "Sheet1"

# Column A: Society Name
$ws.Range("A2").Value = "FLASCO (Florida Society of Clinical Oncology)"
$ws.Range("A3").Value = "GASCO (Georgia Society of Clinical Oncology)"
$ws.Range("A4").Value = "IOS (Indiana Oncology Society)"
$ws.Range("A5").Value = "IOWA Oncology Society"
$ws.Range("A6").Value = "MOASC (Medical Oncology Association of Southern California)"

# Column B: Membership count
$ws.Range("B2").Value = 700
$ws.Range("B3").Value = 250
$ws.Range("B4").Value = 125
$ws.Range("B5").Value = 170
$ws.Range("B6").Value = 300
$ws.Range("B7").Value = 100

# Column C: Encompasses community sites
$ws.Range("C2").Value = "Yes, FLASCO does encompasses community sites, FLASCO's membership primarily consists of academic and hospital-based oncologists and hematologists."
$ws.Range("C3").Value = "Yes, GASCO encompasses community sites. Community oncology practices are included in GASCO to ensure all aspects of cancer care are represented."
$ws.Range("C4").Value = "No, IOS primarily focuses on academic and research institutions ,although it may collaborate with community sites on specific projects or initiatives."
$ws.Range("C5").Value = "No, IOWA Oncology Society focuses on academic centers, as mentioned on their official website."
$ws.Range("C6").Value = "No, MOASC does not encompass community sites, as it focuses on medical oncology practices and providers in Southern California."
$ws.Range("C7").Value = "Yes, many community sites fall within the scope of the term ""encompasses."""

# Column D: Influential on state/local policy
$ws.Range("D2").Value = "Yes, FLASCO is a professional organization focused on clinical oncology and does play a significant role in shaping state or local policy."
$ws.Range("D3").Value = "No, GASCO is not influential on state or local policy. GASCO is a professional organization focused on supporting oncology professionals and providing education and resources, rather than lobbying for policy changes."
$ws.Range("D4").Value = "No, IOS is primarily focused on advancing education, research, and the practice of oncology in Indiana, rather than on influencing policy."
$ws.Range("D5").Value = "No, lack of public information or evidence of direct policy influence."
$ws.Range("D6").Value = "No, MOASC is a professional organization focused on education and advocacy for medical oncologists in Southern California, not a policy-making entity."
$ws.Range("D7").Value = "No, local policies are typically influenced by local factors such as community needs and demographics rather than on a national scale."

# Column E: Engagement opportunity with leadership
$ws.Range("E2").Value = "Yes, FLASCO provides engagement opportunities with leadership, as they offer various networking events and educational programs where members can interact with the leadership team."
$ws.Range("E3").Value = "Yes, GASCO provides engagement opportunity with leadership through mentorship programs, networking events, and involvement in committees and task forces."
$ws.Range("E4").Value = "Yes, IOS provides engagement opportunities with leadership. This is because the organization values input and involvement from its members in decision-making processes and strategic planning."
$ws.Range("E5").Value = "Yes, the IOWA Oncology Society provides engagement opportunities with leadership through networking events, conferences, and mentorship programs."
$ws.Range("E6").Value = "Yes, MOASC does provide engagement opportunities with leadership. The association offers various events, meetings, and forums where members can interact with and learn from industry leaders in the field of oncology."
$ws.Range("E7").Value = "Yes, providing engagement opportunity with leadership is important for fostering a strong organizational culture and promoting employee development."

# Column F: Support for clinical trial recruitment
$ws.Range("F2").Value = "Yes, FLASCO does provides support for clinical trial recruitment. FLASCO is an organization that focuses on education, advocacy, and networking for oncology professionals in Florida."
$ws.Range("F3").Value = "No, GASCO does not provide support for clinical trial recruitment, as they primarily focus on education, advocacy, and networking opportunities for oncology professionals."
$ws.Range("F4").Value = "Yes, IOS offers educational resources and conferences on clinical trials, but does not directly provide recruitment support."
$ws.Range("F5").Value = "No, The IOWA Oncology Society does not provide support for clinical trial recruitment."
$ws.Range("F6").Value = "No, MOASC does not provide support for clinical trial recruitment. , MOASC primarily focuses on education, advocacy, and networking within the field of medical oncology in Southern California."
$ws.Range("F7").Value = "Yes,  ClinicalTrials.gov does provide support for clinical trial recruitment by allowing researchers to register their trials and make them visible to participants searching for studies to participate in."

# Column G: Engagement opportunity with payors
$ws.Range("G2").Value = "Yes, FLASCO provides engagement opportunities with payors. FLASCO works closely with payors to advocate for advancements in oncology care and ensure that patients have access to cutting-edge treatments."
$ws.Range("G3").Value = "No, GASCO does not provide engagement opportunities with payors. GASCO is a professional organization for oncologists and does not focus on payor engagement."
$ws.Range("G4").Value = "Yes, IOS provides engagement opportunities with payors, as part of their mission to improve oncology care and communication with key stakeholders."
$ws.Range("G5").Value = "No, they do not. The IOWA Oncology Society is a professional organization focused on promoting high-quality oncology care in Iowa and providing education and networking opportunities for their members. They do not directly engage with payors."
$ws.Range("G6").Value = "No, MOASC does not directly provide engagement opportunities with payors. The organization primarily focuses on education, advocacy, and networking for medical oncologists in Southern California."
$ws.Range("G7").Value = "Yes, providing engagement opportunities with payors allows for better communication and understanding of payment processes."

# Column H: Area experts on board
$ws.Range("H2").Value = "Yes, The FLASCO board includes area experts from various oncology fields who bring their expertise and knowledge to the organization."
$ws.Range("H3").Value = "Yes, GASCO includes area experts on its board. GASCO is a professional organization for clinical oncologists in Georgia, so it would make sense for the board to include experts in the field."
$ws.Range("H4").Value = "Yes, `n`nThe Indiana Oncology Society includes area experts on its board to ensure that decisions and advancements within the field of oncology are informed and guided by reputable professionals."
$ws.Range("H5").Value = "No, the IOWA Oncology Society does not include area experts on its board. The organization primarily consists of oncologists and healthcare professionals related to oncology."
$ws.Range("H6").Value = "Yes, `nMOASC includes area experts on its board. This can be evidenced by their leadership roles within the organization and their reputation in the field of medical oncology."
$ws.Range("H7").Value = "No, The board of does not include area experts because it consists mainly of industry professionals."

# Column I: Therapeutic research collaborations
$ws.Range("I2").Value = "Yes, FLASCO is involved in therapeutic research collaborations. FLASCO partners with academic institutions, pharmaceutical companies, and other organizations to conduct clinical trials and research studies to improve cancer treatment options."
$ws.Range("I3").Value = "No, GASCO is primarily focused on education and advocacy for oncology professionals in Georgia."
$ws.Range("I4").Value = "No, IOS does not have involvement in therapeutic research collaborations. The society focuses on education, advocacy, and networking for oncology professionals in Indiana."
$ws.Range("I5").Value = "Yes, The IOWA Oncology Society is involved in therapeutic research collaborations. They work with various organizations and institutions to advance cancer treatment options."
$ws.Range("I6").Value = "Yes, MOASC is involved in therapeutic research collaborations. This can be seen by their participation in clinical trials, research studies, and collaborations with other healthcare institutions."
$ws.Range("I7").Value = "No, they are focused on individual research projects within the organization."

# Column J: Top therapeutic area experts on board
$ws.Range("J2").Value = "Yes, The FLASCO board does includes top therapeutic area experts, most of them are practicing oncologists from various specialties."
$ws.Range("J3").Value = "No, GASCO does not include top therapeutic area experts on its board. The organization focuses on representing clinical oncologists in Georgia, rather than experts from specific therapeutic areas."
$ws.Range("J4").Value = "No, The Indiana Oncology Society does not include top therapeutic area experts on its board.Members of the Indiana Oncology Society board are healthcare professionals specializing in hematology/oncology, but their expertise may not be specifically focused on individual therapeutic areas."
$ws.Range("J5").Value = "Yes, the IOWA Oncology Society includes top therapeutic area experts on its board. The organization's board members are often experienced and knowledgeable in various oncology-related fields."
$ws.Range("J6").Value = "No, MOASC does not include top therapeutic area experts on its board. The organization mainly focuses on providing networking and educational opportunities for oncologists in Southern California."
$ws.Range("J7").Value = "No, the board does not include top therapeutic area experts. The board members may have expertise in other areas or fields, but do not specifically specialize in therapeutic areas."

# Column K: Region
$ws.Range("K2").Value = "Florida"
$ws.Range("K3").Value = "Georgia"
$ws.Range("K4").Value = "Midwest"
$ws.Range("K5").Value = "Midwest"
$ws.Range("K6").Value = "Southern California."
$ws.Range("K7").Value = "South America"

